$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark (it currently sits right after
#    "... 0, KEY_SIZE-1, key)" near the top of the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Locate the "4, 4, 4)//rBinarySearch(sortedArray, first, mid-1, key)"
#    phrase further down the document (the second worked example).
$rng = $d.Content
$searchPhrase = "4, 4, 4)//rBinarySearch(sortedArray, first, mid-1, key)"
$found = $rng.Find.Execute($searchPhrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate target phrase"
}
$start = $rng.Start
$end = $rng.End

# 3. Recolor everything from "//" through the end of the phrase green
#    (00B050), matching the comment styling used elsewhere in the doc.
$commentRange = $d.Range($start + 8, $end)
$commentRange.Font.Color = 5287936

# 4. Re-insert the "_GoBack" bookmark so it spans from right after
#    "4, 4, 4)" through the end of the phrase (mirroring where Word
#    last left the edit cursor).
$bmRange = $d.Range($start + 8, $end)
$d.Bookmarks.Add("_GoBack", $bmRange)
